$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:D4").NumberFormat = "@"

$ws.Range("A3").Value = "香"
$ws.Range("B3").Value = "2000"
$ws.Range("C3").Value = "梅煜"
$ws.Range("D3").Value = "2023-10-09"

$ws.Range("A4").Value = "666"
$ws.Range("B4").Value = "0.1"
$ws.Range("C4").Value = "my"
$ws.Range("D4").Value = "2023-10-09"

$ws.Range("A3:D4").ClearFormats()
